$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.670520067214966
$ws.Range("B1").Value = 1.413123250007629
$ws.Range("C1").Value = 1.807768225669861
$ws.Range("D1").Value = 2.7420814037323
$ws.Range("E1").Value = 5.137547492980957
